$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so numeric-looking strings (e.g. "1.007") are not
# auto-converted to numbers by Excel, matching the original inline-string cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.291.16"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.684.82"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "218.04"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "0.5258"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "0.2711"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").Value = "0.06417"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").Value = "0.07488"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "1.723.96"
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("D13").Value = "4.555"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "0.5805"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "0.000008471"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "64.23"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "26.322.51"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "4.927"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "10.86"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "189.04"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "6.199"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "144.51"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "7.713"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("E26").Value = "  +4.60%  "
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").Value = "0.06652"
$ws.Range("E28").Value = "  +12.81%  "
$ws.Range("E29").Value = "  +5.81%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "3.576"
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").Value = "3.569"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "1.664"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "0.6204"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").Value = "2.705"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").Value = "6.391"
$ws.Range("E38").Value = "  +5.77%  "
$ws.Range("D39").Value = "0.01623"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "1.105.99"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("D41").Value = "0.8779"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "1.014"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "100.65"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "1.832.66"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").Value = "56.76"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").Value = "8.173"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "0.05270"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").Value = "0.4302"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "6.036"
$ws.Range("E51").Value = "  +2.49%  "

# Restore the default "General" number format (original cells had no custom format).
$ws.Range("D2:E51").NumberFormat = "General"
